$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.293.45"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.663.90"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.83%  "
$ws.Range("D5").Value = "'219.10"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "'0.5332"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").Value = "'0.2650"
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").Value = "'0.06422"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").Value = "'20.60"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'0.07832"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "'4.571"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "1.664.66"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "1.893.05"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "'0.5521"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "0.0₅8212"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'65.67"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").Value = "'4.692"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "'193.82"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").Value = "'6.042"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").Value = "'145.89"
$ws.Range("E24").Value = "  +2.89%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "'7.200"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "'16.12"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").Value = "'1.482"
$ws.Range("E28").Value = "  +3.80%  "
$ws.Range("D29").Value = "'0.05850"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "'1.282"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").Value = "'3.622"
$ws.Range("E31").Value = "  +2.98%  "
$ws.Range("D32").Value = "'3.283"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "'1.608"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("D34").Value = "'0.9648"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").Value = "'2.825"
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "'0.5804"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("D38").Value = "'0.01608"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D39").Value = "'0.8663"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("D40").Value = "'5.887"
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").Value = "1.050.27"
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").Value = "'104.63"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("D44").Value = "1.803.55"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").Value = "'57.70"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "'1.014"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈104"
$ws.Range("E47").Value = "  -6.68%  "
$ws.Range("D48").Value = "'0.4382"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "'8.058"
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("D50").Value = "'0.05165"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "'1.415"
$ws.Range("E51").Value = "  -4.28%  "
